## Commit: "Added Samples and Files Tab to all tests"
##
## The "startup" sheet drives a Neo4j/Excel data-refresh harness. Row 4
## (FilesTab) holds the Cypher query (column B, header "query") used to
## populate the Files tab. This updates that query so it:
##   - also matches the file's parent node ("-->(parent)")
##   - formats file size into a human-readable value+unit (Bytes/KB/MB/GB/TB)
##   - renames the "Format" column to "File Format"
##   - returns DISTINCT rows ordered by file name
##   - drops the raw `File ID` column from the result

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newFilesQuery = @'
MATCH (f:file)-->(parent)
MATCH (f)-[:file_of_sample]->(samp)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (d)-[:diagnosis_of_study_subject]->(ss)
MATCH (tp)-[:tp_of_diagnosis]->(d)
WHERE tp.chemotherapy_regimen IN ["Dose dense AC (2 week cycles)"]
WITH
        f, parent,p, ss, d,tp, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent,p, ss, d,tp, s, samp,
        f.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
WITH
        f, parent,p, ss, d,tp, s, samp, unit,
        round(factor * value)/factor AS size
RETURN Distinct
    f.file_name AS `File Name`,
    head(labels(samp)) AS `Association`,
    f.file_description AS `Description`,
    f.file_format AS `File Format`,
     CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    p.program_acronym AS `Program Code`,
    s.study_acronym AS `Arm`,
    ss.study_subject_id AS `Case ID`,
    samp.sample_id AS `Sample ID`
    order by f.file_name
'@

# Row 4 / column B = the "query" cell for the FilesTab row.
$ws.Range("B4").Value = $newFilesQuery

# Mirror the author's final cursor position (row 4, the FilesTab query cell).
$ws.Range("B4").Select() | Out-Null
